$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.774.82'
$ws.Range('E2').Value = '  -3.31%  '
$ws.Range('D3').Value = '2.612.27'
$ws.Range('E3').Value = '  -2.08%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'574.31"
$ws.Range('E5').Value = '  -4.41%  '
$ws.Range('D6').Value = "'156.55"
$ws.Range('E6').Value = '  -2.86%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = "'0.625"
$ws.Range('E8').Value = '  -3.00%  '
$ws.Range('D9').Value = '2.609.46'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('E10').Value = '  -7.05%  '
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('E12').Value = '  -5.24%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('D15').Value = '3.082.25'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('E16').Value = '  -7.71%  '
$ws.Range('D17').Value = '63.635.52'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = '2.604.27'
$ws.Range('E18').Value = '  -1.48%  '
$ws.Range('E19').Value = '  -5.07%  '
$ws.Range('D20').Value = "'7.55"
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = "'4.53"
$ws.Range('E21').Value = '  -6.27%  '
$ws.Range('D22').Value = "'343.64"
$ws.Range('E22').Value = '  -3.79%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = "'67.38"
$ws.Range('E24').Value = '  -4.01%  '
$ws.Range('D25').Value = "'1.83"
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('E26').Value = '  -4.87%  '
$ws.Range('D27').Value = "'597.00"
$ws.Range('E27').Value = '  +2.72%  '
$ws.Range('E28').Value = '  -6.49%  '
$ws.Range('D29').Value = "'1.57"
$ws.Range('E29').Value = '  -3.58%  '
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = "'7.92"
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('E33').Value = '  -4.58%  '
$ws.Range('E34').Value = '  -5.49%  '
$ws.Range('D35').Value = "'6.61"
$ws.Range('E35').Value = '  -2.28%  '
$ws.Range('D36').Value = "'5.42"
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('E37').Value = '  -5.31%  '
$ws.Range('D38').Value = "'19.77"
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').Value = "'154.48"
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  -5.12%  '
$ws.Range('D43').Value = "'2.55"
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('D44').Value = "'41.55"
$ws.Range('E44').Value = '  -3.31%  '
$ws.Range('D45').Value = "'157.43"
$ws.Range('E45').Value = '  -3.22%  '
$ws.Range('D46').Value = "'23.84"
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').Value = "'3.90"
$ws.Range('E47').Value = '  -5.48%  '
$ws.Range('E48').Value = '  -5.05%  '
$ws.Range('E49').Value = '  -2.62%  '
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('E51').Value = '  -4.92%  '
